$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "GST Report"

# Drop the frozen header pane (target sheetView has no <pane>)
$excel.ActiveWindow.FreezePanes = $false

# Remove the TOTAL row, the blank spacer row and the footer rows,
# shrinking the used range down to the header + single data row.
$ws.Range("A3:A6").EntireRow.Delete()

# Strip the custom header/border/fill formatting so the remaining two
# rows fall back to the workbook's default style.
$ws.Range("A1:K2").ClearFormats()
$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 15

Write-Output "done"
